$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Atoms")

# Rename the "TEXT" attribute-type cells for the 'name' attribute to "Name"
$ws.Range("B2").Value = "Name"
$ws.Range("B9").Value = "Name"

# Move the isa-relation target values out of column C into column H
$ws.Range("H4").Value = $ws.Range("C4").Value2
$ws.Range("H5").Value = $ws.Range("C5").Value2
$ws.Range("H6").Value = $ws.Range("C6").Value2
$ws.Range("C4").Clear()
$ws.Range("C5").Clear()
$ws.Range("C6").Clear()

# Update the selection shown in the sheet view
$ws.Activate()
$ws.Range("B3:B6").Select()
